$p = $ppt.ActivePresentation
$hm = $p.HandoutMaster
$shp = $hm.Shapes.Item(2)
$shp.TextFrame.TextRange.Text = "12/13/17"
Write-Output ("New text: " + $shp.TextFrame.TextRange.Text)
